$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-by-row updates of Price (D) and Volume(1h) (E) columns.
# For D-column values that parse as plain numbers, force text typing via
# a temporary "@" (text) number format, then ClearFormats() afterwards so
# no residual style index is left on the cell (cells originally carry no
# explicit style), while the text value itself is preserved.

$ws.Range("D2").Value = "62.944.54"
$ws.Range("E2").Value = "  +2.12%  "
$ws.Range("D3").Value = "3.035.43"
$ws.Range("E3").Value = "  +1.12%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.71"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.54"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +6.29%  "
$ws.Range("D8").Value = "3.029.95"
$ws.Range("E8").Value = "  +0.99%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.519"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.36"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +7.67%  "
$ws.Range("E11").Value = "  +2.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.465"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000235"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.45"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +3.21%  "
$ws.Range("E15").Value = "  +1.92%  "
$ws.Range("D16").Value = "3.537.08"
$ws.Range("E16").Value = "  +1.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.08"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.86%  "
$ws.Range("D18").Value = "62.894.89"
$ws.Range("E18").Value = "  +2.21%  "
$ws.Range("D19").Value = "3.036.82"
$ws.Range("E19").Value = "  +1.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "453.14"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.60%  "
$ws.Range("E21").Value = "  +1.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.694"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.68%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.49"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.20"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.93%  "
$ws.Range("E25").Value = "  +5.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.08"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +6.09%  "
$ws.Range("E27").Value = "  +0.31%  "
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.58"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +7.88%  "
$ws.Range("E30").Value = "  +0.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.22"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +8.08%  "
$ws.Range("E32").Value = "  +0.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.61"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.09%  "
$ws.Range("E34").Value = "  +3.02%  "
$ws.Range("D35").Value = "0.0₃0872"
$ws.Range("E35").Value = "  +6.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.04"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +2.01%  "
$ws.Range("E37").Value = "  +2.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.18"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +10.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.10"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.78%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.69"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.72%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.09"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.127"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +3.64%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.311"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +16.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.40"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +7.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "396.76"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.98%  "
$ws.Range("E46").Value = "  +1.80%  "
$ws.Range("D47").Value = "2.741.32"
$ws.Range("E47").Value = "  +0.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.68"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.39%  "
$ws.Range("E49").Value = "  +0.08%  "
$ws.Range("E50").Value = "  +3.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.30"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +3.33%  "
